$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D32").Value = "[Jenkins] 파이프라인 빌드시 git branch목록을 선택하여 배포!"
$ws.Range("E32").Value = "https://dodonam.tistory.com/427"

$ws.Range("D44").Value = "VC 이직 및 생각 정리"
$ws.Range("E44").Value = "https://engineering-ladder.tistory.com/129"

$ws.Range("D51").Value = "[AWS] boto3 SageMaker에 배포된 모델에 추론 요청 날리는 방법"
$ws.Range("E51").Value = "https://bskyvision.com/entry/AWS-boto3-SageMaker%EC%97%90-%EB%B0%B0%ED%8F%AC%EB%90%9C-%EB%AA%A8%EB%8D%B8%EC%97%90-%EC%B6%94%EB%A1%A0-%EC%9A%94%EC%B2%AD-%EB%82%A0%EB%A6%AC%EB%8A%94-%EB%B0%A9%EB%B2%95"
